# Insert a new data row at row 270 (pushing the existing rows 270-340 down to 271-341)
# and populate it with the new record described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 270:340 down by one row, creating a blank row 270
$ws.Rows("270:270").Insert()

# Populate the newly inserted row with its values
$ws.Range("A270").Value = 4
$ws.Range("B270").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C270").Value = "Los Lagos"
$ws.Range("D270").Value = 45135
$ws.Range("E270").Value = 10
$ws.Range("F270").Value = 100112009
$ws.Range("G270").Value = "Acelga"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 60
$ws.Range("K270").Value = 10000
$ws.Range("L270").Value = 10000
$ws.Range("M270").Value = 10000
$ws.Range("N270").Value = "$/docena de atados (12 kilos)"
$ws.Range("O270").Value = "Región de La Araucanía"
$ws.Range("P270").Value = 833
$ws.Range("Q270").Value = 12
$ws.Range("R270").Value = "Hortaliza"

# Ensure the date cell keeps the same date/time number format as the rest of column D
$ws.Range("D270").NumberFormat = $ws.Range("D271").NumberFormat
